$d = $word.ActiveDocument

# 1. "Declaration" + " of interests" -> "Declaration of interests"
$d.Content.Find.Execute("Declaration of interests", $true, $false, $false, $false, $false, $true, 1, $false, "Declaration of interests", 2)

# 2. "Code availability" split across multiple runs with proofErr spell-check markers
$d.Content.Find.Execute("Code availability", $true, $false, $false, $false, $false, $true, 1, $false, "Code availability", 2)

# 3. "During the preparation..." split across runs
$d.Content.Find.Execute("During the preparation of this work the authors used ChatGPT and Claude Code", $true, $false, $false, $false, $false, $true, 1, $false, "During the preparation of this work the authors used ChatGPT and Claude Code", 2)
